$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '34.590.58'
$cell.ClearFormats()
$ws.Range("E2").Value = '  +2.64%  '
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '1.788.15'
$cell.ClearFormats()
$ws.Range("E3").Value = '  +0.84%  '
$ws.Range("E4").Value = '  +0.19%  '
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '223.05'
$cell.ClearFormats()
$ws.Range("E5").Value = '  -0.39%  '
$ws.Range("E6").Value = '  -0.77%  '
$ws.Range("E7").Value = '  +0.38%  '
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '32.03'
$cell.ClearFormats()
$ws.Range("E8").Value = '  +6.79%  '
$ws.Range("E9").Value = '  +0.92%  '
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '0.0683'
$cell.ClearFormats()
$ws.Range("E10").Value = '  +3.51%  '
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '0.0933'
$cell.ClearFormats()
$ws.Range("E11").Value = '  +1.29%  '
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '2.044.07'
$cell.ClearFormats()
$ws.Range("E12").Value = '  +0.89%  '
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '1.781.80'
$cell.ClearFormats()
$ws.Range("E13").Value = '  +0.62%  '
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '10.88'
$cell.ClearFormats()
$ws.Range("E14").Value = '  +8.04%  '
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '34.605.22'
$cell.ClearFormats()
$ws.Range("E15").Value = '  +2.90%  '
$ws.Range("E16").Value = '  +0.99%  '
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '4.28'
$cell.ClearFormats()
$ws.Range("E17").Value = '  +2.71%  '
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '68.24'
$cell.ClearFormats()
$ws.Range("E18").Value = '  -0.09%  '
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '252.42'
$cell.ClearFormats()
$ws.Range("E19").Value = '  +1.30%  '
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '0.0₃0786'
$cell.ClearFormats()
$ws.Range("E20").Value = '  +7.14%  '
$ws.Range("E21").Value = '  -0.11%  '
$ws.Range("E22").Value = '  +1.19%  '
$ws.Range("E23").Value = '  -0.17%  '
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = '2.13'
$cell.ClearFormats()
$ws.Range("E24").Value = '  +0.02%  '
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '158.80'
$cell.ClearFormats()
$ws.Range("E25").Value = '  +0.30%  '
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '16.31'
$cell.ClearFormats()
$ws.Range("E26").Value = '  -0.51%  '
$ws.Range("E27").Value = '  +1.63%  '
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '0.113'
$cell.ClearFormats()
$ws.Range("E28").Value = '  -0.46%  '
$ws.Range("E29").Value = '  -0.01%  '
$ws.Range("E30").Value = '  +0.15%  '
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '3.73'
$cell.ClearFormats()
$ws.Range("E31").Value = '  -1.55%  '
$ws.Range("E32").Value = '  +0.10%  '
$ws.Range("E33").Value = '  +0.27%  '
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '1.87'
$cell.ClearFormats()
$ws.Range("E34").Value = '  +2.77%  '
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '1.425.91'
$cell.ClearFormats()
$ws.Range("E35").Value = '  -3.72%  '
$ws.Range("E36").Value = '  -1.51%  '
$ws.Range("E37").Value = '  +2.57%  '
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '0.628'
$cell.ClearFormats()
$ws.Range("E38").Value = '  +0.20%  '
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '82.90'
$cell.ClearFormats()
$ws.Range("E39").Value = '  -0.12%  '
$ws.Range("E40").Value = '  +4.16%  '
$ws.Range("E41").Value = '  +0.20%  '
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '0.899'
$cell.ClearFormats()
$ws.Range("E42").Value = '  +1.76%  '
$ws.Range("E43").Value = '  -0.47%  '
$ws.Range("E44").Value = '  -1.18%  '
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '0.0499'
$cell.ClearFormats()
$ws.Range("E45").Value = '  -2.59%  '
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '5.91'
$cell.ClearFormats()
$ws.Range("E46").Value = '  +3.74%  '
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '1.942.38'
$cell.ClearFormats()
$ws.Range("E47").Value = '  +1.23%  '
$ws.Range("E48").Value = '  +0.02%  '
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '103.25'
$cell.ClearFormats()
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '11.88'
$cell.ClearFormats()
$ws.Range("E50").Value = '  +2.00%  '
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '49.45'
$cell.ClearFormats()
$ws.Range("E51").Value = '  -3.26%  '
